$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RFID tag / asset ID values (replace placeholder/test values with real ones,
# skipping IDs that already exist as per the import de-duplication fix)
$ws.Range("A2").Value = "E6394850"
$ws.Range("A3").Value = "E0000000"
$ws.Range("A4").Value = "E2354657"
$ws.Range("B4").Value = "m3m3tag"

# Reflect the last-used cell selection at save time
$ws.Range("B4").Select()
